# Update cryptocurrency price/volume data to the new snapshot values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.618.74"
$ws.Range("E2").Value = "  +0.92%  "
$ws.Range("D3").Value = "1.895.70"
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.65"
$ws.Range("E5").Value = "  -3.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.693"
$ws.Range("E6").Value = "  -5.14%  "
$ws.Range("E7").Value = "  -0.79%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.11"
$ws.Range("E8").Value = "  +8.46%  "
$ws.Range("E9").Value = "  -4.28%  "
$ws.Range("E10").Value = "  -2.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0969"
$ws.Range("E11").Value = "  -1.99%  "
$ws.Range("E12").Value = "  +0.99%  "
$ws.Range("D13").Value = "2.171.50"
$ws.Range("E13").Value = "  -0.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.734"
$ws.Range("E14").Value = "  +1.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.98"
$ws.Range("E15").Value = "  +0.20%  "
$ws.Range("D16").Value = "1.904.79"
$ws.Range("E16").Value = "  -0.18%  "
$ws.Range("D17").Value = "35.621.93"
$ws.Range("E17").Value = "  +0.95%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "73.88"
$ws.Range("E18").Value = "  -1.42%  "
$ws.Range("D19").Value = "0.0₃0827"
$ws.Range("E19").Value = "  -2.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "247.45"
$ws.Range("E20").Value = "  +1.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.91"
$ws.Range("E21").Value = "  -1.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.97"
$ws.Range("E22").Value = "  -2.77%  "
$ws.Range("E23").Value = "  -0.81%  "
$ws.Range("E24").Value = "  +4.66%  "
$ws.Range("E25").Value = "  -10.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.58"
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.54"
$ws.Range("E27").Value = "  -1.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.44"
$ws.Range("E28").Value = "  -2.09%  "
$ws.Range("E29").Value = "  -4.08%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.76"
$ws.Range("E31").Value = "  +6.81%  "
$ws.Range("E32").Value = "  -2.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0583"
$ws.Range("E33").Value = "  -0.99%  "
$ws.Range("E34").Value = "  -0.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.854"
$ws.Range("E36").Value = "  -6.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.02"
$ws.Range("E37").Value = "  -0.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.57"
$ws.Range("E38").Value = "  -21.19%  "
$ws.Range("E39").Value = "  +6.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.24"
$ws.Range("E40").Value = "  +1.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "97.93"
$ws.Range("E41").Value = "  +0.91%  "
$ws.Range("E42").Value = "  -1.50%  "
$ws.Range("E43").Value = "  -2.56%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.38"
$ws.Range("E44").Value = "  -2.72%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "1.297.01"
$ws.Range("E45").Value = "  -2.77%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0818"
$ws.Range("E46").Value = "  +8.40%  "
$ws.Range("E47").Value = "  -1.28%  "
$ws.Range("E48").Value = "  -0.25%  "
$ws.Range("E49").Value = "  +3.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.34"
$ws.Range("E51").Value = "  -5.64%  "
